$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerHoje = "14-04-2023"
$headerOperadora = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$prazo = "10 dias úteis"
$respondido = "NO"
$natureza = "Assistencial"
$opcoes = "Responder  Detalhes"

$rows = @(
    @{ Row = 2; Notificacao = "13/04/2023  08:23:06"; Demanda = 12168095; Protocolo = 8606114; Beneficiario = "MÔNICA ALVES GOMES" },
    @{ Row = 3; Notificacao = "13/04/2023  11:42:27"; Demanda = 12168687; Protocolo = 8606848; Beneficiario = "MARCIO CANDIDO DE OLIVEIRA" },
    @{ Row = 4; Notificacao = "13/04/2023  13:26:27"; Demanda = 12169086; Protocolo = 8607367; Beneficiario = "BARBARA KELLY CARNEIRO LEÃO RODRIGUES" },
    @{ Row = 5; Notificacao = "13/04/2023  15:12:39"; Demanda = 12169516; Protocolo = 8607900; Beneficiario = "ANALIS SOARES SILVA" },
    @{ Row = 6; Notificacao = "13/04/2023  16:38:38"; Demanda = 12169806; Protocolo = 8608265; Beneficiario = "MAICKSON CAIQUE VENANCIO" },
    @{ Row = 7; Notificacao = "13/04/2023  16:46:57"; Demanda = 12169830; Protocolo = 8608288; Beneficiario = "EMILLE FERNANDES CORREA" },
    @{ Row = 8; Notificacao = "13/04/2023  16:58:57"; Demanda = 12169862; Protocolo = 8608371; Beneficiario = "MILENA FREIRE TRAVASSOS COUSSEIRO" },
    @{ Row = 9; Notificacao = "13/04/2023  18:17:42"; Demanda = 12170051; Protocolo = 8608573; Beneficiario = "ROGERIA DORALICE SOARES DA SILVA" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $headerHoje
    $ws.Cells.Item($row, 2).Value = $headerOperadora
    $ws.Cells.Item($row, 3).Value = $r.Notificacao
    $ws.Cells.Item($row, 4).Value = $r.Demanda
    $ws.Cells.Item($row, 5).Value = $r.Protocolo
    $ws.Cells.Item($row, 6).Value = $r.Beneficiario
    $ws.Cells.Item($row, 7).Value = $prazo
    $ws.Cells.Item($row, 8).Value = $respondido
    $ws.Cells.Item($row, 9).Value = $natureza
    $ws.Cells.Item($row, 10).Value = $opcoes
}
